$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "Another comment. Because you're worth it"
$ws.Range("B4").Value = "I-Nixon"

$ws.Range("B5").Select()
